# Update localization-status.xlsx for the new handoff cycle:
# - file id a12b13b4-9421-4aef-950e-f1ddf13f842b -> f3621c04-1fb4-4a02-9775-7e3ee7e7defb
# - new xlf hashes for zh-cn / de-de
# - refreshed handoff timestamps

$wb = $excel.ActiveWorkbook

$oldId = "a12b13b4-9421-4aef-950e-f1ddf13f842b"
$newId = "f3621c04-1fb4-4a02-9775-7e3ee7e7defb"

$newZhXlf = "$newId.4453b1d97e2c9b9d5257939f998595323278fe14.zh-cn.xlf"
$newDeXlf = "$newId.4453b1d97e2c9b9d5257939f998595323278fe14.de-de.xlf"

# The original hyperlink targets (relationship URLs) are not part of this
# edit - only the displayed text / cell values change - so reuse the
# existing (old-id) URL when re-creating each hyperlink.
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25cb0d5a3818bb3b785970a50da6c1d42867d504/e2e/$oldId.md"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkUrl, "", "", "e2e\$newId.md")
$wsOverview.Range("G2").Value = "2016-09-04 01:03:33"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newId.md"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkUrl, "", "", "$newId.md")
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-09-04 01:03:28"

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newId.md"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkUrl, "", "", "$newId.md")
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-09-04 01:03:33"
